$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the static "TimeTaken in Hours" value with a formula computed from minutes
$ws.Range("C2").Formula = "=B2/60"
